$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 22:50"

# --- Update Cataluna row (row 5) totals ---
$ws.Range("B5").Value = 18773
$ws.Range("C5").Value = 4966
$ws.Range("D5").Value = 12135
$ws.Range("E5").Value = 1672

# --- Reorder the Canary Islands block above Asturias (rows 19-25) ---
# Before: row19=Asturias, row20=Gran Canaria, row21=La Palma, row22=Lanzarote,
#         row23=Fuerteventura, row24=La Gomera, row25=El Hierro
# After:  row19=Gran Canaria, row20=La Palma, row21=Lanzarote, row22=Fuerteventura,
#         row23=La Gomera, row24=El Hierro, row25=Asturias
$ws.Range("A19").Value = "Gran Canaria"
$ws.Range("B19").Value = 1262
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = 320
$ws.Range("E19").Value = 11

$ws.Range("A20").Value = "La Palma"
$ws.Range("B20").Value = 1262
$ws.Range("C20").Value = 32
$ws.Range("D20").Value = 47
$ws.Range("E20").Value = 2

$ws.Range("A21").Value = "Lanzarote"
$ws.Range("B21").Value = 1262
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = 3

$ws.Range("A22").Value = "Fuerteventura"
$ws.Range("B22").Value = 1262
$ws.Range("C22").Value = 32
$ws.Range("D22").Value = 31
$ws.Range("E22").Value = 0

$ws.Range("A23").Value = "La Gomera"
$ws.Range("B23").Value = 1262
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = 0

$ws.Range("A24").Value = "El Hierro"
$ws.Range("B24").Value = 1262
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 0

$ws.Range("A25").Value = "Asturias"
$ws.Range("B25").Value = 1236
$ws.Range("C25").Value = 90
$ws.Range("D25").Value = 1091
$ws.Range("E25").Value = 55
